# "Generate Report for Handoff" — refresh the handoff report timestamps and
# bump the priority of the still-pending files from "low" to "ht" now that
# they have been handed off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for the files that were just
# generated (rows 4-7) moves from 16:40:53 to 16:41:22.
$overview.Range("G4:G7").Value = "2016-08-14 16:41:22"

# zh-cn: priority for those same files goes from "low" to "ht" (handed off),
# and their "Latest Handoff Datetime" is refreshed.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-14 16:41:14"

# de-de: priority for those same files also goes from "low" to "ht"; its
# "Latest Handoff Datetime" column happens to mirror the Overview generate
# date, so refresh it to match.
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-14 16:41:22"
